# Rename the "population size" feature label to clarify the unit
# (тыс) -> (тыс. чел.), and update the sheet's selection/scroll state.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Численность населения - popsize (тыс. чел.)"

$ws.Range("D18").Select()
